# fix: renomeando tabela no excel
$wb = $excel.ActiveWorkbook

# Rename sheet "ArquivoData" -> "TemplateInfo"
$wsTemplateInfo = $wb.Worksheets.Item("ArquivoData")
$wsTemplateInfo.Name = "TemplateInfo"

# Sheet "Arquivos": update counters
$wsArquivos = $wb.Worksheets.Item("Arquivos")
$wsArquivos.Range("B2").Value = 289
$wsArquivos.Range("C2").Value = 308

# Sheet "TemplateInfo" (formerly "ArquivoData"): relabel row + update values
$wsTemplateInfo.Range("A2").Value = "quantidade"
$wsTemplateInfo.Range("B2").Value = 60
$wsTemplateInfo.Range("C2").Value = 73
$wsTemplateInfo.Range("D2").Value = 68
$wsTemplateInfo.Range("B3").Value = 201
$wsTemplateInfo.Range("C3").Value = 201
$wsTemplateInfo.Range("D3").Value = 201

# Sheet "TemplatesMes": relabel header + update values
$wsTemplatesMes = $wb.Worksheets.Item("TemplatesMes")
$wsTemplatesMes.Range("B1").Value = "mes"
$wsTemplatesMes.Range("C4").Value = 20
$wsTemplatesMes.Range("D4").Value = 23
$wsTemplatesMes.Range("C5").Value = 37
$wsTemplatesMes.Range("D5").Value = 25
